# Update the LR-pairs table (Myoc -> Fzd3) with the newly recomputed
# TPM-based NATMI statistics. The set of sending/target cluster pairs
# changed ("Resolving-Mac" moves from a sending cluster to a target
# cluster) and every numeric column (E:T) is refreshed with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
  @(2, "ECs", "Myoc", "Fzd3", "ECs", 1, 0.3333333333333333, 0.2439023333333333, 0.731707, 0.006368708152767561, 0.006368708152767561, 2, 0.6666666666666666, 0.3251496666666667, 0.975449, 0.07121046526627427, 0.07121046526627427, 0.07930476238255556, 0.713742861443, 0.0004535186707036922, 0.0004535186707036922),
  @(3, "ECs", "Myoc", "Fzd3", "FAPs", 1, 0.3333333333333333, 0.2439023333333333, 0.731707, 0.006368708152767561, 0.006368708152767561, 3, 1, 1.007819666666667, 3.023459, 0.2207208394324094, 0.2207208394324094, 0.2458095682792222, 2.212286114513, 0.001405706609578885, 0.001405706609578885),
  @(4, "ECs", "Myoc", "Fzd3", "MuSCs", 1, 0.3333333333333333, 0.2439023333333333, 0.731707, 0.006368708152767561, 0.006368708152767561, 3, 1, 3.226895, 9.680685, 0.7067166842615477, 0.7067166842615475, 0.7870472199216667, 7.083424979295001, 0.004500872308753377, 0.004500872308753377),
  @(5, "ECs", "Myoc", "Fzd3", "Resolving-Mac", 1, 0.3333333333333333, 0.2439023333333333, 0.731707, 0.006368708152767561, 0.006368708152767561, 1, 0.3333333333333333, 0.006173333333333333, 0.01852, 0.001352011039768762, 0.001352011039768762, 0.001505690404444444, 0.01355121364, [double]"8.610563731607063E-06", [double]"8.610563731607062E-06"),
  @(6, "FAPs", "Myoc", "Fzd3", "ECs", 3, 1, 37.428193, 112.284579, 0.9773142989029397, 0.9773142989029399, 2, 0.6666666666666666, 0.3251496666666667, 0.975449, 0.07121046526627427, 0.07121046526627427, 12.16976447788567, 109.527880300971, 0.06959500593626097, 0.06959500593626099),
  @(7, "FAPs", "Myoc", "Fzd3", "FAPs", 3, 1, 37.428193, 112.284579, 0.9773142989029397, 0.9773142989029399, 3, 1, 1.007819666666667, 3.023459, 0.2207208394324094, 0.2207208394324094, 37.72086899319567, 339.487820938761, 0.2157136324431535, 0.2157136324431535),
  @(8, "FAPs", "Myoc", "Fzd3", "MuSCs", 3, 1, 37.428193, 112.284579, 0.9773142989029397, 0.9773142989029399, 3, 1, 3.226895, 9.680685, 0.7067166842615477, 0.7067166842615475, 120.776848850735, 1086.991639656615, 0.6906843208020846, 0.6906843208020846),
  @(9, "FAPs", "Myoc", "Fzd3", "Resolving-Mac", 3, 1, 37.428193, 112.284579, 0.9773142989029397, 0.9773142989029399, 1, 0.3333333333333333, 0.006173333333333333, 0.01852, 0.001352011039768762, 0.001352011039768762, 0.2310567114533333, 2.07951040308, 0.001321339721440642, 0.001321339721440642),
  @(10, "MuSCs", "Myoc", "Fzd3", "ECs", 3, 1, 0.6248916666666666, 1.874675, 0.01631699294429263, 0.01631699294429263, 2, 0.6666666666666666, 0.3251496666666667, 0.975449, 0.07121046526627427, 0.07121046526627427, 0.2031833171194444, 1.828649854075, 0.001161940659309593, 0.001161940659309593),
  @(11, "MuSCs", "Myoc", "Fzd3", "FAPs", 3, 1, 0.6248916666666666, 1.874675, 0.01631699294429263, 0.01631699294429263, 3, 1, 1.007819666666667, 3.023459, 0.2207208394324094, 0.2207208394324094, 0.6297781112027777, 5.668003000824999, 0.00360150037967697, 0.00360150037967697),
  @(12, "MuSCs", "Myoc", "Fzd3", "MuSCs", 3, 1, 0.6248916666666666, 1.874675, 0.01631699294429263, 0.01631699294429263, 3, 1, 3.226895, 9.680685, 0.7067166842615477, 0.7067166842615475, 2.016459794708334, 18.148138152375, 0.01153149115070956, 0.01153149115070955),
  @(13, "MuSCs", "Myoc", "Fzd3", "Resolving-Mac", 3, 1, 0.6248916666666666, 1.874675, 0.01631699294429263, 0.01631699294429263, 1, 0.3333333333333333, 0.006173333333333333, 0.01852, 0.001352011039768762, 0.001352011039768762, 0.003857664555555555, 0.034718981, [double]"2.206075459651263E-05", [double]"2.206075459651263E-05")
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $value = $row[$i + 1]
        $ws.Range("$col$r").Value = $value
    }
}

Write-Output "Updated $($rows.Length) data rows for Myoc-Fzd3 LR pairs."
